# Adds 20 new data rows (1193-1212) to sheet1 of the Aragon hospitals
# coronavirus dataset, matching the commit "update Spain provinces data
# and charts 2020.06.01".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=1193; Template=2;   Vals=@(43984, 'Hospital Universitario Miguel Servet', 8, 5, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Aragón Hoy')}
    @{Row=1194; Template=199; Vals=@(43984, 'Hospital Clínico Universitario', 22, $null, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Aragón Hoy')}
    @{Row=1195; Template=2;   Vals=@(43984, 'Hospital Royo Villanova', 6, $null, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Aragón Hoy')}
    @{Row=1196; Template=2;   Vals=@(43984, 'Hospital Nuestra Señora de Gracia', 1, $null, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Aragón Hoy')}
    @{Row=1197; Template=202; Vals=@(43984, 'Hospital General de la Defensa', 2, $null, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Aragón Hoy')}
    @{Row=1198; Template=2;   Vals=@(43984, 'Hospital Obispo Polanco', 5, 1, 'Teruel', 'Teruel', 44216, 'Fuente Aragón Hoy')}
    @{Row=1199; Template=202; Vals=@(43984, 'Hospital de Alcañiz', 3, $null, 'Alcañiz', 'Teruel', 44013, 'Fuente Aragón Hoy')}
    @{Row=1200; Template=2;   Vals=@(43984, 'Hospital de Barbastro', 10, 1, 'Barbastro', 'Huesca', 22048, 'Fuente Aragón Hoy')}
    @{Row=1201; Template=202; Vals=@(43984, 'Hospital San Jorge', 7, 2, 'Huesca', 'Huesca', 22125, 'Fuente Aragón Hoy')}
    @{Row=1202; Template=3;   Vals=@(43984, 'Hospital Sagrado Corazón', $null, $null, 'Huesca', 'Huesca', 22125, 'Fuente Aragón Hoy')}
    @{Row=1203; Template=2;   Vals=@(43984, 'Hospital Ernest Lluch', 1, $null, 'Calatayud', 'Zaragoza', 50067, 'Fuente Aragón Hoy')}
    @{Row=1204; Template=2;   Vals=@(43984, 'Hospital San José', 8, $null, 'Teruel', 'Teruel', 44216, 'Fuente Aragón Hoy')}
    @{Row=1205; Template=202; Vals=@(43984, 'Hospital Ejea – Cinco Villas', $null, $null, 'Ejea de los Caballeros', 'Zaragoza', 50095, 'Fuente Aragón Hoy')}
    @{Row=1206; Template=66;  Vals=@(43984, 'MAZ', $null, $null, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Aragón Hoy')}
    @{Row=1207; Template=67;  Vals=@(43984, 'Hospital Viamed Montecanal', $null, $null, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Aragón Hoy')}
    @{Row=1208; Template=66;  Vals=@(43984, 'Clínica Montpellier', 1, $null, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Aragón Hoy')}
    @{Row=1209; Template=67;  Vals=@(43984, 'Hospital Quirón', 1, $null, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Aragón Hoy')}
    @{Row=1210; Template=66;  Vals=@(43984, 'Hospital San Juan de Dios de Zaragoza', $null, $null, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Aragón Hoy')}
    @{Row=1211; Template=67;  Vals=@(43984, 'Clínica Viamed Santiago', $null, $null, 'Huesca', 'Huesca', 22125, 'Fuente Aragón Hoy')}
    @{Row=1212; Template=202; Vals=@(43984, 'Clínica El Pilar', $null, $null, 'Zaragoza', 'Zaragoza', 50297, 'Fuente Aragón Hoy')}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $tpl = $entry.Template

    # Copy the formatting (styles/fills/number formats) of an existing row
    # that uses the identical style pattern, so the new row matches the
    # file's established banding/format conventions exactly.
    $src = $ws.Range("A" + $tpl + ":H" + $tpl)
    $dst = $ws.Range("A" + $r + ":H" + $r)
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats

    $vals = $entry.Vals
    for ($i = 0; $i -lt 8; $i++) {
        $v = $vals[$i]
        if ($null -ne $v) {
            $ws.Cells.Item($r, $i + 1).Value = $v
        }
    }
}

$excel.CutCopyMode = 0

Write-Host $ws.UsedRange.Address()
